# Generate Report for Handback
# Refreshes the handback-status report: two source-file identifiers change
# (one UUID is replaced outright, the other becomes the new "fffff..." id),
# the xlf payload hash is regenerated (now identical for both rows), and
# every handoff/handback timestamp advances a couple of minutes.
#
# The hyperlink *targets* (the external github blob URLs) are untouched by
# this edit - only the visible display text (and therefore the underlying
# cell text) changes - so the original target URLs are reused verbatim.

$wb = $excel.ActiveWorkbook

# ---- new display text for the renamed files --------------------------------
$new3759Md  = "80ee934a-1725-4691-b7d6-a95df99edd55.md"
$newE637Md  = "fffff93752f3-c1d4-4a58-ba81-835ad360fd2a.md"
$newZhCnXlf = "80ee934a-1725-4691-b7d6-a95df99edd55.a35e7655c5a3689f1b68a5cdd07eb06d0c07d945.zh-cn.xlf"
$newDeDeXlf = "80ee934a-1725-4691-b7d6-a95df99edd55.a35e7655c5a3689f1b68a5cdd07eb06d0c07d945.de-de.xlf"

# ---- new timestamps ----------------------------------------------------------
$zhHandoffTime  = "2016-03-23 21:15:17"
$zhHandbackTime = "2016-03-23 21:15:41"
$deHandoffTime  = "2016-03-23 21:15:22"
$deHandbackTime = "2016-03-23 21:15:48"

# =============================================================================
# Sheet "Overview": A2/A3 are hyperlinked .md file names
# =============================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$ovA2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/7de4d6da4fd467c7948367a677effb5120370602/e2e/3759fa03-66d1-4f79-bc1b-0222080b6e36.md"
$ovA3Addr = "https://github.com/OpenLocalizationTest/oltest/blob/7de4d6da4fd467c7948367a677effb5120370602/e2e/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $ovA2Addr, "", "", $new3759Md) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $ovA3Addr, "", "", $newE637Md) | Out-Null

# =============================================================================
# Detail sheet "zh-cn": A/D/F/G columns are hyperlinked, E/H are plain
# timestamp text.
# =============================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhA2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/7de4d6da4fd467c7948367a677effb5120370602/e2e/3759fa03-66d1-4f79-bc1b-0222080b6e36.md"
$zhD2Addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/306edc3e97d0a0030aa96b8e9af734be032fddd7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3759fa03-66d1-4f79-bc1b-0222080b6e36.9c75d3cbd50720f90f5e90d51ec136e188998ef7.zh-cn.xlf"
$zhF2Addr = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3f72717f826651e69f9442389f47a9e06c7b4737/e2e/3759fa03-66d1-4f79-bc1b-0222080b6e36.md"
$zhG2Addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/59e1ba0d7626b585e9dcf10c2aab79c60f5587a8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3759fa03-66d1-4f79-bc1b-0222080b6e36.9c75d3cbd50720f90f5e90d51ec136e188998ef7.zh-cn.xlf"
$zhA3Addr = "https://github.com/OpenLocalizationTest/oltest/blob/7de4d6da4fd467c7948367a677effb5120370602/e2e/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.md"
$zhD3Addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/306edc3e97d0a0030aa96b8e9af734be032fddd7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.f98cfd0985c24a1d224702d9478f1a7fc990b841.zh-cn.xlf"
$zhF3Addr = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3f72717f826651e69f9442389f47a9e06c7b4737/e2e/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.md"
$zhG3Addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/59e1ba0d7626b585e9dcf10c2aab79c60f5587a8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.f98cfd0985c24a1d224702d9478f1a7fc990b841.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Addr, "", "", $new3759Md)  | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhD2Addr, "", "", $newZhCnXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhF2Addr, "", "", $new3759Md)  | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhG2Addr, "", "", $newZhCnXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Addr, "", "", $newE637Md)  | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhD3Addr, "", "", $newZhCnXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhF3Addr, "", "", $newE637Md)  | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhG3Addr, "", "", $newZhCnXlf) | Out-Null

$wsZh.Range("E2").Value = $zhHandoffTime
$wsZh.Range("H2").Value = $zhHandbackTime
$wsZh.Range("E3").Value = $zhHandoffTime
$wsZh.Range("H3").Value = $zhHandbackTime

# =============================================================================
# Detail sheet "de-de": A/D/F/G columns are hyperlinked, E/H are plain
# timestamp text.
# =============================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$deA2Addr = "https://github.com/OpenLocalizationTest/oltest/blob/7de4d6da4fd467c7948367a677effb5120370602/e2e/3759fa03-66d1-4f79-bc1b-0222080b6e36.md"
$deD2Addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4d5f6d5b8e28875a881f69581aa6e86d9de843f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3759fa03-66d1-4f79-bc1b-0222080b6e36.9c75d3cbd50720f90f5e90d51ec136e188998ef7.de-de.xlf"
$deF2Addr = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/22f6e636880537a37a57ebed9fb5b8e98ec0c498/e2e/3759fa03-66d1-4f79-bc1b-0222080b6e36.md"
$deG2Addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a76a1cab92a17688487b321ac075e84837b8b2b3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3759fa03-66d1-4f79-bc1b-0222080b6e36.9c75d3cbd50720f90f5e90d51ec136e188998ef7.de-de.xlf"
$deA3Addr = "https://github.com/OpenLocalizationTest/oltest/blob/7de4d6da4fd467c7948367a677effb5120370602/e2e/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.md"
$deD3Addr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4d5f6d5b8e28875a881f69581aa6e86d9de843f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.f98cfd0985c24a1d224702d9478f1a7fc990b841.de-de.xlf"
$deF3Addr = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/22f6e636880537a37a57ebed9fb5b8e98ec0c498/e2e/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.md"
$deG3Addr = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a76a1cab92a17688487b321ac075e84837b8b2b3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e6375f0c-5214-44ae-ae4b-69c2ac8cd9d4.f98cfd0985c24a1d224702d9478f1a7fc990b841.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Addr, "", "", $new3759Md)  | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deD2Addr, "", "", $newDeDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deF2Addr, "", "", $new3759Md)  | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deG2Addr, "", "", $newDeDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Addr, "", "", $newE637Md)  | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deD3Addr, "", "", $newDeDeXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deF3Addr, "", "", $newE637Md)  | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deG3Addr, "", "", $newDeDeXlf) | Out-Null

$wsDe.Range("E2").Value = $deHandoffTime
$wsDe.Range("H2").Value = $deHandbackTime
$wsDe.Range("E3").Value = $deHandoffTime
$wsDe.Range("H3").Value = $deHandbackTime
